# Add an "index" loop variable to the datalist row-repeat template:
#  - the "[row:list datalist as data]" tag becomes "[row:list datalist as data, index]"
#  - the per-row "${data.id}" placeholder becomes "${index}"
# The active selection is also moved from the old tag row (A6:B6) to the
# row that now holds the index placeholder (A7:B7).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A6").Value = "[row:list datalist as data, index]"
$ws.Range("A7").Value = '${index}'

$ws.Range("A7:B7").Select()
